# Update countries & provincias Spain
# - Refresh stats for India, Sudan, Taiwan and Mongolia
# - Re-sort the country table (rows 4:218) by "Casos totales" (column B) descending,
#   which is how the sheet is maintained; the stat refresh above changes the rank of
#   Sudan (now above Consejo Danes/Burkina Faso/Uruguay/Mayotte) and Mongolia (now
#   above Puerto Rico).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# India (row 18)
$ws.Range("B18").Value = 42533
$ws.Range("C18").Value = 28
$ws.Range("E18").Value = 29367

# Sudan (row 110)
$ws.Range("B110").Value = 678
$ws.Range("C110").Value = 86
$ws.Range("D110").Value = 61
$ws.Range("E110").Value = 576
$ws.Range("H110").Value = 41

# Taiwan (row 121)
$ws.Range("B121").Value = 436
$ws.Range("C121").Value = 4
$ws.Range("E121").Value = 98

# Mongolia (row 175)
$ws.Range("B175").Value = 40
$ws.Range("C175").Value = 1
$ws.Range("D175").Value = 12
$ws.Range("E175").Value = 28
$ws.Range("H175").Value = 0

# Re-sort the whole country table by "Casos totales" (column B) descending so the
# rows reflect the new ranking produced by the updated figures above.
$dataRange = $ws.Range("A4:H218")
$sortKey = $ws.Range("B4:B218")
$dataRange.Sort($sortKey, 2)
